$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.172.43"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "1.835.80"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6674"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07393"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.76%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2951"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07716"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.816.41"
$ws.Range("E12").Value = "  -1.62%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.018"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.08%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6761"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.196"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.51%  "

$ws.Range("D17").Value = "28.919.19"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008238"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9988"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.262"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.708"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.40%  "

$ws.Range("E26").Value = "  -4.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.506"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.206"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.080"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.195"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05360"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.64%  "

$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7500"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.679"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "1.326.32"
$ws.Range("E37").Value = "  +1.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01803"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.740"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9235"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.971"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.82%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.90%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.08210"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +18.21%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000123"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.961.51"
$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5171"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.758"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.286"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05942"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
